$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the Bortle-class description text in column C (rows 2-10) ---
# The underlying text of each description cell is rewritten to be a direct
# quotation (wrapped in quotes) with a trailing citation to Bortle 2001.
# The row -> description mapping (i.e. which class gets which paragraph)
# is unchanged; only the literal text of the shared string changes.

$ws.Range("C2").Value = "`"The zodiacal light, gegenschein, and zodiacal band  are all visible — the zodiacal light to a striking degree, and the zodiacal band spanning the entire sky. Even with direct vision, the Triangulum Galaxy, M33, is is an obvious naked-eye object. The Scorpius and Sagittarius region of the Milky Way casts obvious diffuse shadows on the ground`" (Bortle 2001, p. 01) ."

$ws.Range("C3").Value = "`"Airglow may be weakly apparent along the horizon. The Triangulum Galaxy, M33, is easily seen with direct vision. The summer Milky Way is highly structured to the unaided eye, and its brightest parts look like veined marble when viewed with ordinary binoculars`" (Bortle 2001, p. 01)."

$ws.Range("C4").Value = "`"Some indication of light pollution is evident along the horizon. Clouds may appear faintly illuminated in the brightest parts of the sky near the horizon but are dark overhead. The Milky Way still appears complex`" (Bortle 2001, p. 01)."

$ws.Range("C5").Value = "`"Fairly obvious light-pollution domes are apparent over population centres in several directions. The zodiacal light is clearly evident but doesn't even extend halfway to the zenith at the beginning or end of twilight. The Milky Way well above the horizon is still impressive but lacks all but the most obvious structure`" (Bortle 2001, p. 01)."

$ws.Range("C6").Value = "`"Only hints of the zodiacal light are seen on the best spring and autumn nights. The Milky Way is very weak or invisible near the horizon and looks rather washed out overhead. Light sources are evident in most if not all directions`" (Bortle 2001, p. 02)."

$ws.Range("C7").Value = "`"No trace of the zodiacal light can be seen, even on the best nights. Any indications of the Milky Way are apparent only toward the zenith`" (Bortle 2001, p. 02)."

$ws.Range("C8").Value = "`"The entire sky background has a vague, greyish white hue. Strong light sources are evident in all directions. The Milky Way is totally invisible or nearly so`" (Bortle 2001, p. 02)."

$ws.Range("C9").Value = "`"The sky glows whitish Gray or orangish, and you can read newspaper headlines without difficulty. Some of the stars making up the familiar constellation patterns are difficult to see or are absent entirely`" (Bortle 2001, p. 02)."

$ws.Range("C10").Value = "`"The entire sky is brightly lit, even at the zenith. Many stars making up familiar constellation figures are invisible, and dim constellations such as Cancer and Pisces are not seen at all. Aside from perhaps the Pleiades, no Messier objects are visible to the unaided eye. The only celestial objects that really provide pleasing telescopic views are the Moon, the planets, and a few of the brightest star clusters (if you can find them)`" (Bortle 2001, p. 02)."

# --- Give column C an explicit custom width ---
$ws.Columns.Item(3).ColumnWidth = 8

# --- Move the active selection from G19 to C11 ---
$ws.Range("C11").Select()
